$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Thema" table cell: replace the text content.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(1, 2)
$cellRange = $cell.Range
$cellRange.Find.Execute("Lawinenprobleme im Land Tirol", $true, $false, $false, $false, $false, $true, 1, $false, "Freizeitaktivitäten in Wien", 2)

# ---------------------------------------------------------------------------
# 2) Convert the field-code based HYPERLINK under "Sommer" into a real
#    w:hyperlink run (matching how the other hyperlink in the document,
#    under "folgendem Link abrufbar", is already stored).
# ---------------------------------------------------------------------------
$targetUrl = "https://www.data.gv.at/katalog/dataset/spielplatze-standorte-wien/resource/6f27a91a-bb1e-44f4-a683-98cb80f63379"

$fieldToConvert = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $f = $d.Fields.Item($i)
    if ($f.Code.Text -like "*HYPERLINK*" -and $f.Code.Text -like "*$targetUrl*") {
        $fieldToConvert = $f
    }
}

if ($fieldToConvert -ne $null) {
    $insertStart = $fieldToConvert.Code.Start
    $fieldToConvert.Delete()
    $insertRange = $d.Range($insertStart, $insertStart)
    $d.Hyperlinks.Add($insertRange, $targetUrl, "", "", $targetUrl) | Out-Null
}
